$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the anchor paragraphs by their (stable, pre-edit) text so the script
# does not depend on hard-coded paragraph indices.
# ---------------------------------------------------------------------------

function FindParaIndex($text, $startIdx) {
    for ($i = $startIdx; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Step 1: the "_GoBack" bookmark currently sits alone in an otherwise-empty
# paragraph right after the "search tab" Q&A list. Remove the bookmark so the
# paragraph becomes a plain empty paragraph (the bookmark is re-created later
# in its new location).
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# Step 2: insert a new paragraph right after the "much_answer" #1 question
# paragraph, giving the author ("Tac gia") CSS selector for question #1.
# ---------------------------------------------------------------------------
$idx = FindParaIndex "Câu 1: #div_much_answer > div:nth-child(1) > div.right > div.title_question > a`r" 1
$p = $d.Paragraphs($idx)
$p.Range.InsertParagraphAfter()
$newP = $d.Paragraphs($idx + 1)
$newP.Range.InsertBefore("`tTác giả:  #div_much_answer > div:nth-child(1) > div.right > div.question_info > div.author > strong > a > span")

# ---------------------------------------------------------------------------
# Step 3: after the "need_answer" #20 question paragraph (end of document),
# append: an empty paragraph, a breadcrumb paragraph (carrying the
# "_GoBack" bookmark), and two further empty paragraphs.
# ---------------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$lastP = $d.Paragraphs($lastIdx)
$lastP.Range.InsertParagraphAfter()

$breadcrumbIdx = $lastIdx + 1
$breadcrumbP = $d.Paragraphs($breadcrumbIdx)
$breadcrumbP.Range.InsertBefore("#ctl00_cphMain_ctl00_LeftPane > div.breadcrumb_and_send_question > div.breadcrumb > a:nth-child(3) > span")
$breadcrumbP.Range.InsertParagraphBefore()

# re-resolve after the structural insert shifted indices
$breadcrumbIdx = $breadcrumbIdx + 1
$breadcrumbP = $d.Paragraphs($breadcrumbIdx)
$d.Bookmarks.Add("_GoBack", $breadcrumbP.Range)

$breadcrumbP.Range.InsertParagraphAfter()
$d.Paragraphs($breadcrumbIdx + 1).Range.InsertParagraphAfter()
